$d = $word.ActiveDocument

# 1. Remove the old "_GoBack" bookmark near the top of the document.
#    (Word will automatically renumber the remaining bookmark ids.)
$oldGoBack = $d.Bookmarks("_GoBack")
$oldGoBack.Delete()

# 2. Fix the duplicated "to to request" -> "to request" and split the run so a
#    new "_GoBack" bookmark (marking the last edit point) sits between the two
#    halves of the sentence.
$r = $d.Content
$found = $r.Find.Execute("You DO NOT need to to request", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Position right after "You DO NOT need to " (before the duplicated "to ").
    $splitStart = $r.Start + 19
    # Position right after the duplicated "to " (before "request").
    $dupEnd = $r.Start + 22

    # Insert the bookmark first, at the zero-length split point. This splits the
    # run into two (both still carrying the original run's formatting/rsid)
    # without touching any run text yet.
    $bmRange = $d.Range($splitStart, $splitStart)
    $d.Bookmarks.Add("_GoBack", $bmRange)

    # Now remove the duplicated "to " - it lives entirely inside the new
    # (second) run, so only that run's text is rewritten.
    $dupRange = $d.Range($splitStart, $dupEnd)
    $dupRange.Text = ""
}
